# Nhat Linh : update final flow for report function admin
#
# Updates the "8/2023" (rows 30 & 33) and "10/2023" (rows 38 & 41) revenue
# report blocks on the DuLieuThongKe sheet, and refreshes the best-fit width
# of column I ("Tong thue") to match the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "8/2023" block: Top 1 doanh thu (row 30) ---
$ws.Range("E30").Value = 1351990.0
$ws.Range("F30").Value = 1065000.0
$ws.Range("G30").Value = 151414.29498244822
$ws.Range("H30").Value = 2.0
$ws.Range("I30").Value = 135199.0020146221
$ws.Range("J30").Value = 376.7030029296875

# --- "8/2023" block: Toan cua hang (row 33) ---
$ws.Range("E33").Value = 1351990.0
$ws.Range("F33").Value = 1065000.0
$ws.Range("G33").Value = 151414.29498244822
$ws.Range("H33").Value = 2.0
$ws.Range("I33").Value = 135199.0020146221
$ws.Range("J33").Value = 376.7030029296875

# --- "10/2023" block: Top 1 doanh thu (row 38) ---
$ws.Range("E38").Value = 0.0
$ws.Range("F38").Value = 0.0
$ws.Range("G38").Value = 0.0
$ws.Range("H38").Value = 0.0
$ws.Range("I38").Value = 0.0
$ws.Range("J38").Value = 0.0

# --- "10/2023" block: Toan cua hang (row 41) ---
$ws.Range("E41").Value = 0.0
$ws.Range("F41").Value = 0.0
$ws.Range("G41").Value = 0.0
$ws.Range("H41").Value = 0.0
$ws.Range("I41").Value = 0.0
$ws.Range("J41").Value = 0.0

# Column I's best-fit width shrinks now that its widest value changed.
$ws.Columns.Item(9).ColumnWidth = 11.62890625
